# Add a new "Multiple Choice" style question table to sheet "1_"
# and make that sheet the active/selected tab (mirroring the structure
# already present on sheet "0_").

$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("0_")
$ws1 = $wb.Worksheets.Item("1_")

# --- Populate sheet "1_" with the new question table ----------------------
# (Values are written in the same order the original author typed them —
# column A top-to-bottom, then the "wrong answer" comments in column C,
# then the last row — so newly created shared-string entries land at the
# same indices as the authored workbook.)

$ws1.Range("A1").Value = "Why do you think the coefficient of drag changes with velocity?"
$ws1.Range("A2").Value = "Because the shape of the baseball changes very slightly, and this changes C_d"
$ws1.Range("A3").Value = "Because the shape of the baseball changes very slightly, and this changes the cross-sectional area of the ball"
$ws1.Range("A4").Value = "Because drag depends on the square of velocity, and so as velocity increases, the coefficient has to take this into account"

$ws1.Range("B1").Value = "Correct"
$ws1.Range("B2").Value = "N"
$ws1.Range("B3").Value = "N"
$ws1.Range("B4").Value = "N"

$ws1.Range("C1").Value = "Comment"
$ws1.Range("C2").Value = "This is absolute nonsense!"
$ws1.Range("C3").Value = "This answer is rubbish!"
$ws1.Range("C4").Value = "This makes no sense whatsoever!"

$ws1.Range("A5").Value = "Because the drag equation is just a model that tries to simplify a complex physical system, and the model has limitations"
$ws1.Range("B5").Value = "Y"
$ws1.Range("C5").Value = "Yahoo, Bob!  It's a model!  So many of those equations that you think are *Truth* and *Law* are just models we use to help us understand a complex world!"

# --- Row heights & wrap-text style (matches sheet "0_" layout) ------------

$ws1.Range("A1:C5").WrapText = $true

$ws1.Rows.Item(1).RowHeight = 45
$ws1.Rows.Item(2).RowHeight = 60
$ws1.Rows.Item(3).RowHeight = 75
$ws1.Rows.Item(4).RowHeight = 75
$ws1.Rows.Item(5).RowHeight = 75

# --- Selections / active sheet ---------------------------------------------

$ws0.Range("A1:C5").Select()

$ws1.Activate()
$ws1.Range("C12").Select()
